# --- Replace separate "Kinh do"/"Vi do" (longitude/latitude) columns (Q, R) ---
# with a single combined "Toa do" (coordinates) column (Q) ---
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column R entirely (old "Vi do" / latitude column)
$ws.Range("R:R").Delete()

# Rename header Q1 from "Kinh do" to "Toa do"
$ws.Range("Q1").Value = "Tọa độ"

# Replace the numeric longitude in Q2 with the combined lat/long text
$ws.Range("Q2").ClearFormats()
$ws.Range("Q2").Value = "21.079350776626914, 105.80247286566104"

# Widen column Q to fit the longer coordinate text
$ws.Columns("Q:Q").ColumnWidth = 41.83

# --- Fix customer phone numbers (typos) ---
# Row 2: 0962547000 -> 0962447000
$ws.Range("B2").Value = "'0962447000"
# Row 3: 0962547001 -> 0912547001
$ws.Range("B3").Value = "'0912547001"

# --- Restore view selection state ---
$null = $ws.Range("F9").Select()
